# Apply the "Added calcium and b12 to Zambia and Uganda and updated output"
# re-run metadata refresh to spade_uganda_h_iron.xlsx:
#  - Info sheet: refresh Start_time / End_time stamps
#  - sessionInfo sheet: bump package version numbers (here, magrittr, rprojroot)
#    and drop the "backports" row from the Loaded_only package table

$wb = $excel.ActiveWorkbook

$infoWs = $wb.Worksheets.Item("Info")
$sessionWs = $wb.Worksheets.Item("sessionInfo")

# Info!B26/B27 - Start_time / End_time of the (re-)run
$infoWs.Range("B26").Value = "Thu Nov 19 15:23:47 2020"
$infoWs.Range("B27").Value = "Thu Nov 19 15:23:54 2020"

# sessionInfo - "Ohter_packages" / "Loaded_only" version bumps
$sessionWs.Range("G2").Value = "1.0.0"   # here
$sessionWs.Range("J3").Value = "2.0.1"   # magrittr
$sessionWs.Range("J10").Value = "2.0.2"  # rprojroot

# sessionInfo - remove the "backports" entry (row 15) from the Loaded_only
# table; the row below ("boot" / "1.3-25", previously row 16) shifts up to
# take its place, and the now-trailing row 16 cells are cleared entirely.
$sessionWs.Range("I15").Value = "boot"
$sessionWs.Range("J15").Value = "1.3-25"
$sessionWs.Range("I16").ClearContents()
$sessionWs.Range("J16").ClearContents()
